$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.223.15'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  -2.78%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.872.63'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  -4.12%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '515.68'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -5.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.15'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -8.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.534'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -6.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.880.62'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -4.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.10'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -0.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.105'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -8.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.352'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -4.81%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.380.11'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  -4.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.127'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +1.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.414.63'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  -2.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.06'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -8.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.893.66'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -3.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000138'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -6.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.81'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -7.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.30'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -6.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '346.44'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -8.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.42'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -4.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.15'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -3.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.439'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -6.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.175'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -7.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.01'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +1.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.61'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -8.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0838'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -10.82%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.65'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -4.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.23'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -6.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '153.23'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -4.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.28'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  -7.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.47'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -7.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.966'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -9.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.17'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -8.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.19'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -0.91%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.644'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -4.41%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.253.13'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -6.99%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.41'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -9.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.61'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -7.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0570'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -4.03%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.85'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -10.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.77'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -9.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0232'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -5.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.34'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0898'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -5.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.99'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -8.85%  '
